$wb = $excel.ActiveWorkbook

$wsCustomers = $wb.Worksheets.Item("Customers")
$wsAddresses = $wb.Worksheets.Item("Addresses")

$wsCustomers.Range("A2").Value = "CUST1"
$wsCustomers.Range("B2").Value = "XX"
$wsCustomers.Range("C2").Value = "DT11"

$wsAddresses.Range("A2").Value = "CUST1"
$wsAddresses.Range("B2").Value = "XX"
$wsAddresses.Range("C2").Value = "addr1"

$wsCustomers.Columns.Item(1).ColumnWidth = 13.1
$wsAddresses.Columns.Item(1).ColumnWidth = 13.1

$wsCustomers.Range("C2").Select()
$wsAddresses.Range("D8").Select()
